$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update aggregated values in row 2 (g_H aggregation)
$ws.Range("B2").Value = 566.432
$ws.Range("D2").Value = 566.432
$ws.Range("F2").Value = 34181.24137931035

# Remove the now-redundant row 3 (id_DK_Decentral_HS), folded into row 2's aggregate
$ws.Rows.Item(3).Delete()
